$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Globo"
$ws.Range("B6").Value = "RJ TV 1"
$ws.Range("C6").Value = "Transporte"
$ws.Range("D6").Value = "2025-04-01T11:52"
$ws.Range("E6").Value = "Negativo"
$ws.Range("F6").Value = "Sem vans de novo. Permissionários do setor C, em Campos, fazem paralisação alegando falta de pagamento. Repórter *ao vivo* do ponto das vans. Mais uma paralisação. Muitas pessoas foram surpreendidas ao voltar para casa. Permissionários pedem o repasse feito pela prefeitura e alegam inconsistência da empresa de bilhetagem. Só estão rodando as vans do Setor C, que vão para Travessão. Entrevista com passageiros no ponto, reclamando de falta de transporte e dos custos sem o repasse da prefeitura. Passageiros pagam R$ 3,50, utilizando cartão Rio Card e precisam da prefeitura para complementar. Motoristas foram para frente da prefeitura hoje de manhã. Entrevista com permissionário Jefferson Henrique. Em nota, prefeitura informou que não há atraso de pagamento dos subsídios. Tema foi abordado ontem à noite."
